$d = $word.ActiveDocument

$oldText = "Intressant att få uppleva musikindustrin från insidan. Uppträden i TV-produktioner såsom Allsång på skansen och Melodifestivalen."
$firstNewText = "Uppträden i TV-produktioner såsom Allsång på skansen och Melodifestivalen."
$secondNewText = "Intressant att få uppleva musikindustrin från insidan. "

# Locate the paragraph that currently holds the combined sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ($oldText + "`r")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph containing the source text."
}

# Range covering just the paragraph's visible text (exclude the trailing
# paragraph mark) so we can safely overwrite it in place.
$r = $d.Range($target.Range.Start, $target.Range.End - 1)
$r.Text = $firstNewText

# Insert a brand-new paragraph right after, inheriting the ListBullet style,
# and give it the remaining sentence.
$newPara = $r.InsertParagraphAfter()

$found = $false
$applied = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $p.Range.Text = $secondNewText
        $applied = $true
        break
    }
    if ($p.Range.Text -eq ($firstNewText + "`r")) {
        $found = $true
    }
}

if (-not $applied) {
    throw "Could not locate newly inserted paragraph to set its text."
}
